$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column N (14th column), shifting N:T -> O:U.
$ws.Columns.Item(14).Insert()

# Row 1 (headers): split "NameMesEr" into "FirstNameMesEr" (M1) / "LastNameMesEr" (new N1)
$ws.Range("M1").Value = "FirstNameMesEr"
$ws.Range("N1").Value = "LastNameMesEr"

# Row 2: add validation messages
$ws.Range("N2").Value = "You can't leave this empty."
$ws.Range("O2").Value = "You can't leave this empty."

# Row 3: add validation message
$ws.Range("N3").Value = "You can't leave this empty."

# Row 4: clear stale UserName test value, add validation messages
$ws.Range("D4").ClearContents()
$ws.Range("M4").Value = "You can't leave this empty."
$ws.Range("O4").Value = "You can't leave this empty."

# Row 5: add validation message
$ws.Range("N5").Value = "You can't leave this empty."

# Update the view: scrolled to column M, normal zoom 100%, selection on P13
$ws.Application.ActiveWindow.ScrollColumn = 13
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("P13").Select()
